# Append a new block of paragraphs (notes + a "To Do" bulleted list) to the
# end of the document, right after the "More Visualizations of Input
# Embedding and Latent Space" paragraph and before the section break.

$d = $word.ActiveDocument

# Common run/paragraph-mark formatting used throughout this block: matches
# the formatting already used by the rest of the document's body text.
$rPr = '<w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr>'

# List-paragraph properties for the "To Do" bullets -- reuses the numbering
# definition (numId 1) already present in the document's numbering part.
$listPPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' + $rPr + '</w:pPr>'
$plainPPr = '<w:pPr>' + $rPr + '</w:pPr>'

function New-EmptyPara {
    return '<w:p>' + $plainPPr + '</w:p>'
}

function New-TextPara {
    param([string]$text)
    return '<w:p>' + $plainPPr + '<w:r>' + $rPr + '<w:t>' + $text + '</w:t></w:r></w:p>'
}

function New-ListPara {
    param([string]$text)
    return '<w:p>' + $listPPr + '<w:r>' + $rPr + '<w:t>' + $text + '</w:t></w:r></w:p>'
}

function New-ListParaSpellErr {
    param([string]$text)
    return '<w:p>' + $listPPr + '<w:proofErr w:type="spellStart"/><w:r>' + $rPr + '<w:t>' + $text + '</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
}

$paras = @()
$paras += New-EmptyPara
$paras += New-TextPara "Make it robust so it works with any dataset"
$paras += New-EmptyPara
$paras += New-EmptyPara
$paras += New-EmptyPara
$paras += New-TextPara "To Do"
$paras += New-ListPara "Build Dataset Builder"
$paras += New-ListParaSpellErr "Visualisations"
$paras += New-ListPara "Build Custom dataset"
$paras += New-ListPara "Write Academic Paper"
$paras += New-ListPara "Input Synthesis"
$paras += New-ListPara "Images to Smiles"
$paras += New-ListPara "Full Automatic Pipeline"
$paras += New-ListPara "UI"

$body = [string]::Join("", $paras)

$openXmlPackage = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">' +
'<pkg:xmlData>' +
'<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">' +
'<Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/>' +
'</Relationships>' +
'</pkg:xmlData>' +
'</pkg:part>' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' + $body +
'<w:sectPr><w:pgSz w:w="12240" w:h="15840"/><w:pgMar w:top="1440" w:right="1440" w:bottom="1440" w:left="1440" w:header="720" w:footer="720" w:gutter="0"/><w:cols w:space="720"/></w:sectPr>' +
'</w:body>' +
'</w:document>' +
'</pkg:xmlData>' +
'</pkg:part>' +
'</pkg:package>'

$insertionPoint = $d.Paragraphs.Last.Range
$insertionPoint.Collapse(0)
$insertionPoint.InsertXML($openXmlPackage) | Out-Null

Write-Host "Inserted" $paras.Count "new paragraphs. Document now has" $d.Paragraphs.Count "paragraphs."
